$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 131
$ws.Range("A131").Value = 129
$ws.Range("B131").Value = 7483081
$ws.Range("C131").Value = "Ecuador LigaPro Serie A"
$ws.Range("D131").Value = "Ecuador LigaPro Serie A"
$ws.Range("E131").Value = 45255.83333333334
$ws.Range("F131").Value = "Deportivo Cuenca"
$ws.Range("G131").Value = "El Nacional"
$ws.Range("H131").Value = 1
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = "H"
$ws.Range("K131").Value = 2.75
$ws.Range("L131").Value = 3.25
$ws.Range("M131").Value = 2.55
$ws.Range("N131").Value = 3
$ws.Range("O131").Value = 3.3
$ws.Range("P131").Value = 2.3
$ws.Range("Q131").Value = 0.25
$ws.Range("R131").Value = 1.825
$ws.Range("S131").Value = 1.975
$ws.Range("T131").Value = 2.75
$ws.Range("U131").Value = 2
$ws.Range("V131").Value = 1.8
$ws.Range("W131").Value = 2
$ws.Range("X131").Value = -1
$ws.Range("Y131").Value = -1
$ws.Range("Z131").Value = 0.825
$ws.Range("AA131").Value = -1
$ws.Range("AB131").Value = -1
$ws.Range("AC131").Value = 0.8

# Row 132
$ws.Range("A132").Value = 130
$ws.Range("B132").Value = 7483189
$ws.Range("C132").Value = "Ecuador LigaPro Serie A"
$ws.Range("D132").Value = "Ecuador LigaPro Serie A"
$ws.Range("E132").Value = 45255.83333333334
$ws.Range("F132").Value = "Independiente del Valle"
$ws.Range("G132").Value = "Orense"
$ws.Range("H132").Value = 2
$ws.Range("I132").Value = 2
$ws.Range("J132").Value = "D"
$ws.Range("K132").Value = 1.4
$ws.Range("L132").Value = 4.75
$ws.Range("M132").Value = 7
$ws.Range("N132").Value = 1.4
$ws.Range("O132").Value = 4.5
$ws.Range("P132").Value = 8
$ws.Range("Q132").Value = -1.25
$ws.Range("R132").Value = 1.875
$ws.Range("S132").Value = 1.925
$ws.Range("T132").Value = 2.5
$ws.Range("U132").Value = 1.925
$ws.Range("V132").Value = 1.875
$ws.Range("W132").Value = -1
$ws.Range("X132").Value = 3.5
$ws.Range("Y132").Value = -1
$ws.Range("Z132").Value = -1
$ws.Range("AA132").Value = 0.925
$ws.Range("AB132").Value = 0.925
$ws.Range("AC132").Value = -1

# Row 133
$ws.Range("A133").Value = 131
$ws.Range("B133").Value = 7483281
$ws.Range("C133").Value = "Ecuador LigaPro Serie A"
$ws.Range("D133").Value = "Ecuador LigaPro Serie A"
$ws.Range("E133").Value = 45255.83333333334
$ws.Range("F133").Value = "SD Aucas"
$ws.Range("G133").Value = "Delfin SC"
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = "D"
$ws.Range("K133").Value = 1.909
$ws.Range("L133").Value = 3.25
$ws.Range("M133").Value = 4.2
$ws.Range("N133").Value = 1.909
$ws.Range("O133").Value = 3.5
$ws.Range("P133").Value = 4
$ws.Range("Q133").Value = -0.5
$ws.Range("R133").Value = 1.9
$ws.Range("S133").Value = 1.9
$ws.Range("T133").Value = 2.5
$ws.Range("U133").Value = 1.8
$ws.Range("V133").Value = 2
$ws.Range("W133").Value = -1
$ws.Range("X133").Value = 2.5
$ws.Range("Y133").Value = -1
$ws.Range("Z133").Value = -1
$ws.Range("AA133").Value = 0.8999999999999999
$ws.Range("AB133").Value = -1
$ws.Range("AC133").Value = 1

# Row 136
$ws.Range("A136").Value = 134
$ws.Range("B136").Value = 7482832
$ws.Range("C136").Value = "Ecuador LigaPro Serie A"
$ws.Range("D136").Value = "Ecuador LigaPro Serie A"
$ws.Range("E136").Value = 45256.83333333334
$ws.Range("F136").Value = "Barcelona Guayaquil"
$ws.Range("G136").Value = "Guayaquil City"
$ws.Range("H136").Value = 2
$ws.Range("I136").Value = 1
$ws.Range("J136").Value = "H"
$ws.Range("K136").Value = 1.363
$ws.Range("L136").Value = 5
$ws.Range("M136").Value = 7.5
$ws.Range("N136").Value = 1.444
$ws.Range("O136").Value = 4
$ws.Range("P136").Value = 8
$ws.Range("Q136").Value = -1.25
$ws.Range("R136").Value = 2.05
$ws.Range("S136").Value = 1.75
$ws.Range("T136").Value = 2.5
$ws.Range("U136").Value = 1.95
$ws.Range("V136").Value = 1.85
$ws.Range("W136").Value = 0.444
$ws.Range("X136").Value = -1
$ws.Range("Y136").Value = -1
$ws.Range("Z136").Value = -0.5
$ws.Range("AA136").Value = 0.375
$ws.Range("AB136").Value = 0.95
$ws.Range("AC136").Value = -1

# Row 137
$ws.Range("A137").Value = 135
$ws.Range("B137").Value = 7483306
$ws.Range("C137").Value = "Ecuador LigaPro Serie A"
$ws.Range("D137").Value = "Ecuador LigaPro Serie A"
$ws.Range("E137").Value = 45256.83333333334
$ws.Range("F137").Value = "Tecnico Universitario"
$ws.Range("G137").Value = "Club Atletico Libertad"
$ws.Range("H137").Value = 1
$ws.Range("I137").Value = 1
$ws.Range("J137").Value = "D"
$ws.Range("K137").Value = 1.5
$ws.Range("L137").Value = 4.333
$ws.Range("M137").Value = 5.75
$ws.Range("N137").Value = 1.533
$ws.Range("O137").Value = 4.2
$ws.Range("P137").Value = 5.5
$ws.Range("Q137").Value = -1
$ws.Range("R137").Value = 1.925
$ws.Range("S137").Value = 1.875
$ws.Range("T137").Value = 2.25
$ws.Range("U137").Value = 1.8
$ws.Range("V137").Value = 2
$ws.Range("W137").Value = -1
$ws.Range("X137").Value = 3.2
$ws.Range("Y137").Value = -1
$ws.Range("Z137").Value = -1
$ws.Range("AA137").Value = 0.875
$ws.Range("AB137").Value = -0.5
$ws.Range("AC137").Value = 0.5

# Row 142
$ws.Range("A142").Value = 140
$ws.Range("B142").Value = 7528852
$ws.Range("C142").Value = "Ecuador LigaPro Serie A"
$ws.Range("D142").Value = "Ecuador LigaPro Serie A"
$ws.Range("E142").Value = 45263.83333333334
$ws.Range("F142").Value = "Delfin SC"
$ws.Range("G142").Value = "Tecnico Universitario"
$ws.Range("H142").Value = 2
$ws.Range("I142").Value = 2
$ws.Range("J142").Value = "D"
$ws.Range("K142").Value = 2.1
$ws.Range("L142").Value = 3.4
$ws.Range("M142").Value = 3.1
$ws.Range("N142").Value = 2.1
$ws.Range("O142").Value = 3.4
$ws.Range("P142").Value = 3.1
$ws.Range("Q142").Value = -0.25
$ws.Range("R142").Value = 1.8
$ws.Range("S142").Value = 2
$ws.Range("T142").Value = 2.25
$ws.Range("U142").Value = 1.9
$ws.Range("V142").Value = 1.9
$ws.Range("W142").Value = -1
$ws.Range("X142").Value = 2.4
$ws.Range("Y142").Value = -1
$ws.Range("Z142").Value = -0.5
$ws.Range("AA142").Value = 0.5
$ws.Range("AB142").Value = 0.8999999999999999
$ws.Range("AC142").Value = -1

# Row 143
$ws.Range("A143").Value = 141
$ws.Range("B143").Value = 7528857
$ws.Range("C143").Value = "Ecuador LigaPro Serie A"
$ws.Range("D143").Value = "Ecuador LigaPro Serie A"
$ws.Range("E143").Value = 45263.83333333334
$ws.Range("F143").Value = "Universidad Catolica del Ecuador"
$ws.Range("G143").Value = "Barcelona Guayaquil"
$ws.Range("H143").Value = 0
$ws.Range("I143").Value = 1
$ws.Range("J143").Value = "A"
$ws.Range("K143").Value = 1.533
$ws.Range("L143").Value = 4
$ws.Range("M143").Value = 5.5
$ws.Range("N143").Value = 1.5
$ws.Range("O143").Value = 4.333
$ws.Range("P143").Value = 5.25
$ws.Range("Q143").Value = -1
$ws.Range("R143").Value = 1.8
$ws.Range("S143").Value = 2
$ws.Range("T143").Value = 3
$ws.Range("U143").Value = 1.975
$ws.Range("V143").Value = 1.825
$ws.Range("W143").Value = -1
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = 4.25
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 1
$ws.Range("AB143").Value = -1
$ws.Range("AC143").Value = 0.825

# Row 144
$ws.Range("A144").Value = 142
$ws.Range("B144").Value = 7528858
$ws.Range("C144").Value = "Ecuador LigaPro Serie A"
$ws.Range("D144").Value = "Ecuador LigaPro Serie A"
$ws.Range("E144").Value = 45263.83333333334
$ws.Range("F144").Value = "Orense"
$ws.Range("G144").Value = "SD Aucas"
$ws.Range("H144").Value = 1
$ws.Range("I144").Value = 2
$ws.Range("J144").Value = "A"
$ws.Range("K144").Value = 2.2
$ws.Range("L144").Value = 3.2
$ws.Range("M144").Value = 3.2
$ws.Range("N144").Value = 1.95
$ws.Range("O144").Value = 3.2
$ws.Range("P144").Value = 3.8
$ws.Range("Q144").Value = -0.5
$ws.Range("R144").Value = 1.95
$ws.Range("S144").Value = 1.85
$ws.Range("T144").Value = 2.25
$ws.Range("U144").Value = 1.85
$ws.Range("V144").Value = 1.95
$ws.Range("W144").Value = -1
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = 2.8
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 0.8500000000000001
$ws.Range("AB144").Value = 0.8500000000000001
$ws.Range("AC144").Value = -1

# Row 145
$ws.Range("A145").Value = 143
$ws.Range("B145").Value = 7528848
$ws.Range("C145").Value = "Ecuador LigaPro Serie A"
$ws.Range("D145").Value = "Ecuador LigaPro Serie A"
$ws.Range("E145").Value = 45263.83333333334
$ws.Range("F145").Value = "Emelec"
$ws.Range("G145").Value = "Deportivo Cuenca"
$ws.Range("H145").Value = 2
$ws.Range("I145").Value = 1
$ws.Range("J145").Value = "H"
$ws.Range("K145").Value = 1.75
$ws.Range("L145").Value = 3.5
$ws.Range("M145").Value = 4.2
$ws.Range("N145").Value = 2.4
$ws.Range("O145").Value = 3.1
$ws.Range("P145").Value = 2.75
$ws.Range("Q145").Value = -0.25
$ws.Range("R145").Value = 2.05
$ws.Range("S145").Value = 1.75
$ws.Range("T145").Value = 2.25
$ws.Range("U145").Value = 1.8
$ws.Range("V145").Value = 2
$ws.Range("W145").Value = 1.4
$ws.Range("X145").Value = -1
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 1.05
$ws.Range("AA145").Value = -1
$ws.Range("AB145").Value = 0.8
$ws.Range("AC145").Value = -1

# Row 209
$ws.Range("A209").Value = 207
$ws.Range("B209").Value = 7773510
$ws.Range("C209").Value = "Ecuador LigaPro Serie A"
$ws.Range("D209").Value = "Ecuador LigaPro Serie A"
$ws.Range("E209").Value = 45399.875
$ws.Range("F209").Value = "Tecnico Universitario"
$ws.Range("G209").Value = "Deportivo Cuenca"
$ws.Range("H209").Value = 2
$ws.Range("I209").Value = 1
$ws.Range("J209").Value = "H"
$ws.Range("K209").Value = 1.95
$ws.Range("L209").Value = 3.25
$ws.Range("M209").Value = 4.2
$ws.Range("N209").Value = 2.4
$ws.Range("O209").Value = 3.1
$ws.Range("P209").Value = 3
$ws.Range("Q209").Value = -0.25
$ws.Range("R209").Value = 2.05
$ws.Range("S209").Value = 1.75
$ws.Range("T209").Value = 2.25
$ws.Range("U209").Value = 1.9
$ws.Range("V209").Value = 1.9
$ws.Range("W209").Value = 1.4
$ws.Range("X209").Value = -1
$ws.Range("Y209").Value = -1
$ws.Range("Z209").Value = 1.05
$ws.Range("AA209").Value = -1
$ws.Range("AB209").Value = 0.8999999999999999
$ws.Range("AC209").Value = -1

# Row 210
$ws.Range("A210").Value = 208
$ws.Range("B210").Value = 7773772
$ws.Range("C210").Value = "Ecuador LigaPro Serie A"
$ws.Range("D210").Value = "Ecuador LigaPro Serie A"
$ws.Range("E210").Value = 45400.66666666666
$ws.Range("F210").Value = "Imbabura"
$ws.Range("G210").Value = "Universidad Catolica del Ecuador"
$ws.Range("H210").Value = 2
$ws.Range("I210").Value = 2
$ws.Range("J210").Value = "D"
$ws.Range("K210").Value = 4.5
$ws.Range("L210").Value = 3.75
$ws.Range("M210").Value = 1.727
$ws.Range("N210").Value = 4.5
$ws.Range("O210").Value = 4
$ws.Range("P210").Value = 1.533
$ws.Range("Q210").Value = 1
$ws.Range("R210").Value = 1.85
$ws.Range("S210").Value = 1.95
$ws.Range("T210").Value = 2.75
$ws.Range("U210").Value = 1.9
$ws.Range("V210").Value = 1.9
$ws.Range("W210").Value = -1
$ws.Range("X210").Value = 3
$ws.Range("Y210").Value = -1
$ws.Range("Z210").Value = 0.8500000000000001
$ws.Range("AA210").Value = -1
$ws.Range("AB210").Value = 0.8999999999999999
$ws.Range("AC210").Value = -1

# Row 211
$ws.Range("A211").Value = 209
$ws.Range("B211").Value = 7773068
$ws.Range("C211").Value = "Ecuador LigaPro Serie A"
$ws.Range("D211").Value = "Ecuador LigaPro Serie A"
$ws.Range("E211").Value = 45400.77083333334
$ws.Range("F211").Value = "Independiente del Valle"
$ws.Range("G211").Value = "Mushuc Runa"
$ws.Range("H211").Value = 2
$ws.Range("I211").Value = 1
$ws.Range("J211").Value = "H"
$ws.Range("K211").Value = 1.333
$ws.Range("L211").Value = 5.5
$ws.Range("M211").Value = 8
$ws.Range("N211").Value = 1.4
$ws.Range("O211").Value = 5.25
$ws.Range("P211").Value = 6
$ws.Range("Q211").Value = -1.25
$ws.Range("R211").Value = 1.875
$ws.Range("S211").Value = 1.925
$ws.Range("T211").Value = 3
$ws.Range("U211").Value = 1.8
$ws.Range("V211").Value = 2
$ws.Range("W211").Value = 0.3999999999999999
$ws.Range("X211").Value = -1
$ws.Range("Y211").Value = -1
$ws.Range("Z211").Value = -0.5
$ws.Range("AA211").Value = 0.4625
$ws.Range("AB211").Value = 0
$ws.Range("AC211").Value = -0
